$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (14 and 15) to the feed logs sheet
$newRows = @(
    @(13, 1, "2024-06-14 17:33:02", 200, 0),
    @(14, 2, "2024-06-14 17:33:02", 200, 0)
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
